$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.221.88"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.645.81"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.31"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.95"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.876.52"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.52"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.543"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.66"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.194.05"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.54"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("E22").Value = "  +5.07%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.19"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.92"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.77"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260.33"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.851"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  +6.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.41"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.786.35"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.84"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.01%  "

# Row 50/51 swap: Algorand moves to rank 50, EnergySwap moves to rank 51
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0973"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +1.08%  "
